$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.549.60"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.640.24"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'309.05"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.3769"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'52.75"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").Value = "'0.3684"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").Value = "'1.279"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "'0.08203"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "'0.9988"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'23.20"
$ws.Range("E13").Value = "  +3.99%  "
$ws.Range("D14").Value = "'6.671"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").Value = "'7.467"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "1.640.34"
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").Value = "'95.06"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").Value = "'0.06956"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").Value = "'18.41"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'6.583"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").Value = "'0.9979"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "23.558.68"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").Value = "'12.97"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +7.35%  "
$ws.Range("D26").Value = "'2.409"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "'151.18"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'5.328"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'135.88"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "'2.413"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").Value = "'6.871"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").Value = "1.819.73"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("D34").Value = "'0.9772"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "'0.02818"
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("D36").Value = "'10.48"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("D37").Value = "'0.07472"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'6.237"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("D39").Value = "'0.2543"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'0.08867"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").Value = "'1.402"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("D42").Value = "'0.7175"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "'12.64"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "'16.19"
$ws.Range("E44").Value = "  +8.75%  "
$ws.Range("D45").Value = "'0.6629"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("D46").Value = "'2.361"
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("D47").Value = "'4.051"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D48").Value = "'0.9986"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'0.08065"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").Value = "'131.10"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = "  +0.34%  "
